$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for the columns that move (A,B,D,E,F,G,H,K,Q,R) across rows 5-17
# before overwriting any of them, since this is a permutation of whole-row data.
$snapshot = @{}
foreach ($r in 5..17) {
    $rowData = @{}
    $rowData["A"] = $ws.Cells.Item($r, 1).Value2
    $rowData["B"] = $ws.Cells.Item($r, 2).Value2
    $rowData["D"] = $ws.Cells.Item($r, 4).Value2
    $rowData["E"] = $ws.Cells.Item($r, 5).Value2
    $rowData["F"] = $ws.Cells.Item($r, 6).Value2
    $rowData["G"] = $ws.Cells.Item($r, 7).Value2
    $rowData["H"] = $ws.Cells.Item($r, 8).Value2
    $rowData["K"] = $ws.Cells.Item($r, 11).Value2
    $rowData["Q"] = $ws.Cells.Item($r, 17).Value2
    $rowData["R"] = $ws.Cells.Item($r, 18).Value2
    $snapshot[$r] = $rowData
}

# Mapping: destination row -> source row (data that should end up there)
$mapping = @{
    5 = 7
    6 = 17
    7 = 11
    8 = 9
    9 = 12
    10 = 5
    11 = 16
    12 = 10
    13 = 15
    14 = 6
    15 = 8
    16 = 13
    17 = 14
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $data = $snapshot[$srcRow]
    $ws.Cells.Item($destRow, 1).Value2 = $data["A"]
    $ws.Cells.Item($destRow, 2).Value2 = $data["B"]
    $ws.Cells.Item($destRow, 4).Value2 = $data["D"]
    $ws.Cells.Item($destRow, 5).Value2 = $data["E"]
    $ws.Cells.Item($destRow, 6).Value2 = $data["F"]
    $ws.Cells.Item($destRow, 7).Value2 = $data["G"]
    $ws.Cells.Item($destRow, 8).Value2 = $data["H"]
    if ($data["K"] -eq $null -or $data["K"] -eq "") {
        $ws.Cells.Item($destRow, 11).ClearContents()
    } else {
        $ws.Cells.Item($destRow, 11).Value2 = $data["K"]
    }
    $ws.Cells.Item($destRow, 17).Value2 = $data["Q"]
    $ws.Cells.Item($destRow, 18).Value2 = $data["R"]
}
